$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B3").Value = 1200.03
$ws1.Range("B4").Value = 0.03
$ws1.Range("B5").Value = 0.04
$ws1.Range("B6").Value = 14
$ws1.Range("B7").Value = 5
$ws1.Range("B9").Value = 35.71

# --- Strategy Status sheet ---
$ws2 = $wb.Worksheets.Item("Strategy Status")
$ws2.Range("C4").Value = 100.03
$ws2.Range("D4").Value = 14
$ws2.Range("E4").Value = 0.03
$ws2.Range("F4").Value = 0.03
$ws2.Range("G4").Value = 35.71

# --- All Trades sheet: add row 15 for trade #14 ---
$ws3 = $wb.Worksheets.Item("All Trades")
# Seed row 15 from row 14 so text cells keep matching string/style types
# (direct string assignment for full dates gets auto-parsed into date serials)
$ws3.Range("A14:Q14").Copy($ws3.Range("A15:Q15"))
$ws3.Range("A15").Value = 14
$ws3.Range("C15").Value = "04:07:26"
$ws3.Range("G15").Value = 0.806631
$ws3.Range("I15").Value = 0.8288
$ws3.Range("J15").Value = 0.01
$ws3.Range("K15").Value = 100.03

# --- MarketMaking sheet: add row 15 for trade #14 ---
$ws4 = $wb.Worksheets.Item("MarketMaking")
$ws4.Range("A14:Q14").Copy($ws4.Range("A15:Q15"))
$ws4.Range("A15").Value = 14
$ws4.Range("C15").Value = "04:07:26"
$ws4.Range("G15").Value = 0.806631
$ws4.Range("I15").Value = 0.8288
$ws4.Range("J15").Value = 0.01
$ws4.Range("K15").Value = 100.03
